$wb = $excel.ActiveWorkbook

# --- Sheet: snapshot (sheet1) - update scraped_at timestamps in column K ---
$ws1 = $wb.Worksheets.Item("snapshot")

$timestamps = @{
    2  = "2025-11-29T07:01:55.304932+00:00"
    3  = "2025-11-29T07:01:57.567316+00:00"
    4  = "2025-11-29T07:01:57.567347+00:00"
    5  = "2025-11-29T07:01:57.567366+00:00"
    6  = "2025-11-29T07:01:59.812483+00:00"
    7  = "2025-11-29T07:02:02.593664+00:00"
    8  = "2025-11-29T07:02:05.067873+00:00"
    9  = "2025-11-29T07:02:07.550978+00:00"
    10 = "2025-11-29T07:02:13.190660+00:00"
    11 = "2025-11-29T07:02:13.190692+00:00"
    12 = "2025-11-29T07:02:15.745059+00:00"
    13 = "2025-11-29T07:02:18.075103+00:00"
    14 = "2025-11-29T07:02:18.075135+00:00"
    15 = "2025-11-29T07:02:20.926624+00:00"
    16 = "2025-11-29T07:02:23.198179+00:00"
    17 = "2025-11-29T07:02:23.198211+00:00"
    18 = "2025-11-29T07:02:23.198230+00:00"
    19 = "2025-11-29T07:02:25.491485+00:00"
    20 = "2025-11-29T07:02:25.491519+00:00"
    21 = "2025-11-29T07:02:25.491543+00:00"
    22 = "2025-11-29T07:02:28.123706+00:00"
    23 = "2025-11-29T07:02:28.123737+00:00"
    24 = "2025-11-29T07:02:30.908385+00:00"
    25 = "2025-11-29T07:02:30.908415+00:00"
    26 = "2025-11-29T07:02:30.908433+00:00"
    27 = "2025-11-29T07:02:30.908451+00:00"
    28 = "2025-11-29T07:02:33.193094+00:00"
    29 = "2025-11-29T07:02:38.684107+00:00"
    30 = "2025-11-29T07:02:44.208773+00:00"
    31 = "2025-11-29T07:02:44.208802+00:00"
    32 = "2025-11-29T07:02:46.473504+00:00"
    33 = "2025-11-29T07:02:46.473533+00:00"
}

foreach ($row in $timestamps.Keys) {
    $ws1.Cells.Item($row, 11).Value = $timestamps[$row]
}

# --- Sheet: new_injured (sheet3) - remove data rows 2 and 3 (keep header) ---
$ws3 = $wb.Worksheets.Item("new_injured")
$ws3.Rows.Item(3).Delete()
$ws3.Rows.Item(2).Delete()
